$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp text (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 12 de Abril de 2020 a las 12:52"

# --- Malta / Kirguistan swap (rows 102-103) with updated Malta numbers ---
# Row 102 now shows Malta with refreshed case numbers
$ws.Range("A102").Value = "Malta"
$ws.Range("B102").Value = 378
$ws.Range("C102").Value = 8
$ws.Range("D102").Value = 44
$ws.Range("E102").Value = 331
$ws.Range("F102").Value = 4
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 3

# Row 103 now shows Kirguistan, keeping its previous numbers
$ws.Range("A103").Value = "Kirguistan"
$ws.Range("B103").Value = 377
$ws.Range("C103").Value = 38
$ws.Range("D103").Value = 44
$ws.Range("E103").Value = 328
$ws.Range("F103").Value = 5
$ws.Range("G103").Value = 0
$ws.Range("H103").Value = 5

# --- Laos / Sudan / Angola reorder (rows 169-171) with updated Laos numbers ---
# Row 169 now shows Laos with refreshed case numbers
$ws.Range("A169").Value = "Laos"
$ws.Range("B169").Value = 19
$ws.Range("C169").Value = 1
$ws.Range("D169").Value = 0
$ws.Range("E169").Value = 19
$ws.Range("F169").Value = 0
$ws.Range("G169").Value = 0
$ws.Range("H169").Value = 0

# Row 170 now shows Sudan, keeping its previous numbers
$ws.Range("A170").Value = "Sudan"
$ws.Range("B170").Value = 19
$ws.Range("C170").Value = 0
$ws.Range("D170").Value = 2
$ws.Range("E170").Value = 15
$ws.Range("F170").Value = 0
$ws.Range("G170").Value = 0
$ws.Range("H170").Value = 2

# Row 171 now shows Angola, keeping its previous numbers
$ws.Range("A171").Value = "Angola"
$ws.Range("B171").Value = 19
$ws.Range("C171").Value = 0
$ws.Range("D171").Value = 4
$ws.Range("E171").Value = 13
$ws.Range("F171").Value = 0
$ws.Range("G171").Value = 0
$ws.Range("H171").Value = 2
